# "Generate Report for Handback" - mark the zh-cn / de-de handoff rows as
# handed back: update status, record the handback target/file columns
# (mirroring the original handoff file, since the handback is "in sync"
# with the source), and stamp the handback datetime.

$wb = $excel.ActiveWorkbook

$sheetInfo = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-01-28 04:35:57";
       MdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/f552e4a454097815a4eb83986a6dd40a9353fed6/e2e/5349622b-63c6-4b75-8a9c-d2552ded057f.md";
       XlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4955338b446238dacfe4dd4abf620789f44be2f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/5349622b-63c6-4b75-8a9c-d2552ded057f.203c81c31ece72fc4fb7b650d7f31a24fc5c5ce3.zh-cn.xlf";
       MdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/f552e4a454097815a4eb83986a6dd40a9353fed6/e2e/c75298cf-bbb3-41e9-a4e2-80ef06980076.md";
       XlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4955338b446238dacfe4dd4abf620789f44be2f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/c75298cf-bbb3-41e9-a4e2-80ef06980076.88f3b8499bec165cbf106484ec6597e245dc729d.zh-cn.xlf" },
    @{ Name = "de-de"; HandbackTime = "2016-01-28 04:36:14";
       MdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/f552e4a454097815a4eb83986a6dd40a9353fed6/e2e/5349622b-63c6-4b75-8a9c-d2552ded057f.md";
       XlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e568275785d0fffdede00ce4e90411d3eadd4bd8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/5349622b-63c6-4b75-8a9c-d2552ded057f.203c81c31ece72fc4fb7b650d7f31a24fc5c5ce3.de-de.xlf";
       MdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/f552e4a454097815a4eb83986a6dd40a9353fed6/e2e/c75298cf-bbb3-41e9-a4e2-80ef06980076.md";
       XlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e568275785d0fffdede00ce4e90411d3eadd4bd8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/c75298cf-bbb3-41e9-a4e2-80ef06980076.88f3b8499bec165cbf106484ec6597e245dc729d.de-de.xlf" }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Row 2 ------------------------------------------------------------
    $mdName2  = $ws.Range("A2").Value2
    $xlfName2 = $ws.Range("C2").Value2

    $ws.Range("B2").Value = "Handed back: in sync with en-US"

    $ws.Range("E2").Value = $mdName2
    $ws.Hyperlinks.Add($ws.Range("E2"), $info.MdUrl1, "", "", $mdName2)

    $ws.Range("F2").Value = $xlfName2
    $ws.Hyperlinks.Add($ws.Range("F2"), $info.XlfUrl1, "", "", $xlfName2)

    $ws.Range("G2").Value = $info.HandbackTime

    # Row 3 ------------------------------------------------------------
    $mdName3  = $ws.Range("A3").Value2
    $xlfName3 = $ws.Range("C3").Value2

    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    $ws.Range("E3").Value = $mdName3
    $ws.Hyperlinks.Add($ws.Range("E3"), $info.MdUrl2, "", "", $mdName3)

    $ws.Range("F3").Value = $xlfName3
    $ws.Hyperlinks.Add($ws.Range("F3"), $info.XlfUrl2, "", "", $xlfName3)

    $ws.Range("G3").Value = $info.HandbackTime
}

Write-Output "Handback report generated for zh-cn and de-de sheets."
